$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings si 6 and si 9) ---
$ws.Range("A8").Value = "Volume 30   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/7/2023  Through  8/13/2023"

# --- Cells that change TYPE (number <-> text placeholder) ---
# Copy number-format from a donor cell that already carries the target style,
# then write the value, so the cell lands on the same style index as in the target file.
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").Value = "0"

$ws.Range("G15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1

$ws.Range("D22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 1

$ws.Range("D23").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = "0"

$ws.Range("D27").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 2

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("H14").Value = -100
$ws.Range("I15").Value = 8
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 14.285714285714
$ws.Range("N15").Value = -57.894736842105
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 24
$ws.Range("H16").Value = -29.166666666666
$ws.Range("I16").Value = 124
$ws.Range("J16").Value = 152
$ws.Range("K16").Value = -18.421052631578
$ws.Range("L16").Value = 21.568627450980
$ws.Range("M16").Value = 53.086419753086
$ws.Range("N16").Value = -85.167464114832
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -42.105263157894
$ws.Range("I17").Value = 118
$ws.Range("J17").Value = 112
$ws.Range("K17").Value = 5.357142857142
$ws.Range("L17").Value = 43.902439024390
$ws.Range("M17").Value = 114.545454545455
$ws.Range("N17").Value = -32.183908045977
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -71.428571428571
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -35
$ws.Range("I18").Value = 152
$ws.Range("J18").Value = 151
$ws.Range("K18").Value = 0.662251655629
$ws.Range("L18").Value = 25.619834710743
$ws.Range("M18").Value = 0.662251655629
$ws.Range("N18").Value = -90.759878419452
$ws.Range("C19").Value = 52
$ws.Range("D19").Value = 30
$ws.Range("E19").Value = 73.333333333333
$ws.Range("F19").Value = 156
$ws.Range("G19").Value = 143
$ws.Range("H19").Value = 9.090909090909
$ws.Range("I19").Value = 1037
$ws.Range("J19").Value = 1041
$ws.Range("K19").Value = -0.384245917387
$ws.Range("L19").Value = 65.92
$ws.Range("M19").Value = 36.627140974967
$ws.Range("N19").Value = -54.094732182381
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = -61.111111111111
$ws.Range("I20").Value = 105
$ws.Range("J20").Value = 118
$ws.Range("K20").Value = -11.016949152542
$ws.Range("L20").Value = 7.142857142857
$ws.Range("M20").Value = 110
$ws.Range("N20").Value = -94.902912621359
$ws.Range("C21").Value = 64
$ws.Range("D21").Value = 48
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 205
$ws.Range("G21").Value = 226
$ws.Range("H21").Value = -9.292035398230
$ws.Range("I21").Value = 1547
$ws.Range("J21").Value = 1584
$ws.Range("K21").Value = -2.335858585858
$ws.Range("L21").Value = 49.180327868852
$ws.Range("M21").Value = 40.126811594202
$ws.Range("N21").Value = -77.887364208118
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -33.333333333333
$ws.Range("I22").Value = 32
$ws.Range("J22").Value = 28
$ws.Range("K22").Value = 14.285714285714
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 88.235294117647
$ws.Range("M23").Value = -23.529411764705
$ws.Range("C24").Value = 77
$ws.Range("D24").Value = 83
$ws.Range("E24").Value = -7.228915662650
$ws.Range("F24").Value = 279
$ws.Range("G24").Value = 282
$ws.Range("H24").Value = -1.063829787234
$ws.Range("I24").Value = 2028
$ws.Range("J24").Value = 2371
$ws.Range("K24").Value = -14.466469843947
$ws.Range("L24").Value = 35.109926715523
$ws.Range("M24").Value = 99.018645731108
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 233.333333333333
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 23
$ws.Range("H25").Value = -4.347826086956
$ws.Range("I25").Value = 201
$ws.Range("J25").Value = 229
$ws.Range("K25").Value = -12.227074235807
$ws.Range("L25").Value = 2.551020408163
$ws.Range("M25").Value = -10.666666666666
$ws.Range("I26").Value = 18
$ws.Range("K26").Value = 50
$ws.Range("L26").Value = 38.461538461538
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 62
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = 24
$ws.Range("L27").Value = 10.714285714285
$ws.Range("D30").Value = 1
$ws.Range("G30").Value = 4
$ws.Range("J30").Value = 23
$ws.Range("K30").Value = -73.913043478260
